# 20170601 Route 1 Fixes.
# Each timepoint block of 5 identical rows is trimmed down to 4 rows by
# deleting the last (5th) row of each block (rows 8,13,18,...,113 before
# any shifting takes place). Deleting bottom-up keeps the row numbers of
# not-yet-processed rows stable.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowsToDelete = @(113,108,103,98,93,88,83,78,73,68,63,58,53,48,43,38,33,28,23,18,13,8)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# Excel's row-delete/shift does not perfectly preserve the alternating
# zebra-stripe formatting (cell style index 4) on the shifted rows, so
# reassert the banding explicitly: odd data rows (5,7,9,...) are shaded,
# even data rows (4,6,8,...) are not. Use copy/paste-special of formats
# from two known-good reference rows so no new style entries are created.

$ws.Range("A4:D4").Copy()
for ($r = 4; $r -le 91; $r += 2) {
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)
}

$ws.Range("A5:D5").Copy()
for ($r = 5; $r -le 91; $r += 2) {
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)
}

$excel.CutCopyMode = 0
